# "Задание 1-1" is the active sheet in this workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the formula in G3: the denominator of the division (C5/SIN(C3/C5)*SIN(C3/C5))
# was missing a closing/opening paren pair around the full SIN(...)*SIN(...) product,
# i.e. it should divide by (SIN(C3/C5)*SIN(C3/C5)) as a whole rather than only by
# SIN(C3/C5) and then multiplying by SIN(C3/C5) again.
$ws.Range("G3").Formula = "=C4*C5*C3*C3 - (C5/(SIN(C3/C5)*SIN(C3/C5)))"

# Make column G wide enough to show the (now longer) formula result comfortably.
$ws.Columns.Item(7).ColumnWidth = 56.8

# Move the active selection to G5.
$ws.Range("G5").Select()

# Restore/refresh the workbook window layout.
$win = $wb.Windows.Item(1)
$win.Left = -225
$win.Top = 825
$win.Width = 10890
$win.Height = 9660
